$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "88.649.38"
$ws.Range("E2").Value = "  +10.35%  "
$ws.Range("D3").Value = "3.359.55"
$ws.Range("E3").Value = "  +5.34%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "222.43"
$ws.Range("E5").Value = "  +6.28%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "648.81"
$ws.Range("E6").Value = "  +3.42%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.343"
$ws.Range("E7").Value = "  +25.20%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.618"
$ws.Range("E9").Value = "  +5.11%  "
$ws.Range("D10").Value = "3.361.35"
$ws.Range("E10").Value = "  +5.50%  "
$ws.Range("E11").Value = "  +4.68%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000277"
$ws.Range("E12").Value = "  +7.43%  "
$ws.Range("E13").Value = "  +2.28%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.37"
$ws.Range("E14").Value = "  +11.22%  "
$ws.Range("D15").Value = "3.987.73"
$ws.Range("E15").Value = "  +5.73%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.49"
$ws.Range("E16").Value = "  +3.70%  "
$ws.Range("D17").Value = "88.423.75"
$ws.Range("E17").Value = "  +9.95%  "
$ws.Range("D18").Value = "3.359.93"
$ws.Range("E18").Value = "  +5.47%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.76"
$ws.Range("E19").Value = "  +3.68%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.21"
$ws.Range("E20").Value = "  +6.37%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "470.01"
$ws.Range("E21").Value = "  +6.76%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.30"
$ws.Range("E22").Value = "  +1.30%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.58"
$ws.Range("E23").Value = "  +7.59%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.52"
$ws.Range("E24").Value = "  +24.32%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "7.43"
$ws.Range("E25").Value = "  +7.17%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "5.50"
$ws.Range("E26").Value = "  +16.85%  "
$ws.Range("D27").Value = "3.507.91"
$ws.Range("E27").Value = "  +4.50%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "79.47"
$ws.Range("E28").Value = "  +4.23%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.202"
$ws.Range("E29").Value = "  +63.54%  "
$ws.Range("E30").Value = "  +5.66%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.998"
$ws.Range("E31").Value = "  -0.17%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "9.42"
$ws.Range("E32").Value = "  +5.33%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "597.71"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.997"
$ws.Range("E34").Value = "  -0.13%  "
$ws.Range("E35").Value = "  +7.45%  "
$ws.Range("E36").Value = "  +4.33%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.152"
$ws.Range("E37").Value = "  +0.45%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "24.27"
$ws.Range("E38").Value = "  +5.06%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.92"
$ws.Range("E39").Value = "  +22.49%  "
$ws.Range("E40").Value = "  +4.05%  "
$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.12"
$ws.Range("E41").Value = "  +17.27%  "
$ws.Range("B42").Value = "FirstDigitalUSD"
$ws.Range("C42").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.00"
$ws.Range("E42").Value = "  -0.31%  "
$ws.Range("B43").Value = "WhiteBITCoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "21.71"
$ws.Range("E43").Value = "  +4.58%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.08"
$ws.Range("E44").Value = "  +13.23%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "193.54"
$ws.Range("E45").Value = "  +2.13%  "
$ws.Range("E46").Value = "  +0.04%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "157.69"
$ws.Range("E47").Value = "  -3.20%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "47.61"
$ws.Range("E48").Value = "  +10.99%  "
$ws.Range("E49").Value = "  +7.61%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.799"
$ws.Range("E50").Value = "  +1.64%  "
$ws.Range("B51").Value = "InjectiveProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "26.86"
$ws.Range("E51").Value = "  +6.38%  "
